# Added all tagging and file download functionality
$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Add the new "tags" worksheet right after Sheet1
$tags = $wb.Worksheets.Add($null, $sheet1)
$tags.Name = "tags"

# Header row
$tags.Range("A1").Value = "Statement"
$tags.Range("B1").Value = "Aspect"
$tags.Range("C1").Value = "Sentiment"

# Format header row like Sheet1!A1 (bold, centered, thin box border)
$header = $tags.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2

# Data rows
$tags.Range("A2").Value = "The rooms were bad but staff was friendly."
$tags.Range("B2").Value = "room"
$tags.Range("C2").Value = "NEG"

$tags.Range("A3").Value = "It is hot but scenery is good ."
$tags.Range("B3").Value = "weather"
$tags.Range("C3").Value = "NEG"

$tags.Range("A4").Value = "It is hot but scenery is good ."
$tags.Range("B4").Value = "scenery"
$tags.Range("C4").Value = "POS"

$tags.Range("A5").Value = "The rooms were bad but staff was friendly."
$tags.Range("B5").Value = "staff"
$tags.Range("C5").Value = "POS"

$tags.Activate()
$tags.Range("A1").Select()
